$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 3-7: add column I value (Valor final), mark status Concluído ---
$ws.Range("I3").Value = 772
$ws.Range("K3").Value = "Concluído"

$ws.Range("I4").Value = 761
$ws.Range("K4").Value = "Concluído"

$ws.Range("I5").Value = 751.45
$ws.Range("K5").Value = "Concluído"

$ws.Range("I6").Value = 766.5
$ws.Range("K6").Value = "Concluído"

$ws.Range("I7").Value = 75750
$ws.Range("K7").Value = "Concluído"

# --- Row 8: update Valor final and status ---
$ws.Range("I8").Value = 11355
$ws.Range("K8").Value = "Concluído"

# --- New row 9 ---
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = "fasdfasdf"
$ws.Range("C9").Value = "asdfasdf"
$ws.Range("D9").Value = "00:00"
$ws.Range("E9").Value = 50
$ws.Range("F9").Value = 50
$ws.Range("G9").Value = 50
$ws.Range("H9").Value = 50
$ws.Range("I9").Value = 7575757555
$ws.Range("J9").Value = "dvasdfasdfasdfasdfasdfs"
$ws.Range("K9").Value = "Concluído"

# --- New row 10 ---
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = "sadfasdfasd"
$ws.Range("C10").Value = "'20"
$ws.Range("D10").Value = "'20"
$ws.Range("E10").Value = 20
$ws.Range("F10").Value = 202
$ws.Range("G10").Value = 2
$ws.Range("H10").Value = 2
$ws.Range("I10").Value = 21212102.2
$ws.Range("J10").Value = "02`t02`t02`t02`t02`t02`t02`t0"
$ws.Range("K10").Value = "Concluído"

# --- New row 11 ---
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "asdasdasds"
$ws.Range("C11").Value = "dasdasdasd"
$ws.Range("D11").Value = "'20"
$ws.Range("E11").Value = 2
$ws.Range("F11").Value = 2
$ws.Range("G11").Value = 2
$ws.Range("H11").Value = 20
$ws.Range("I11").Value = 315.7
$ws.Range("K11").Value = "Concluído"

# --- New row 12 ---
$ws.Range("A12").Value = 11
$ws.Range("B12").Value = "asdfasdfasdf"
$ws.Range("C12").Value = "fasdfasdf"
$ws.Range("D12").Value = "'20"
$ws.Range("E12").Value = 20
$ws.Range("F12").Value = 20
$ws.Range("G12").Value = 20
$ws.Range("H12").Value = 20
$ws.Range("I12").Value = 322
$ws.Range("K12").Value = "Concluído"

# --- Update sheet view selection: whole row 9 selected, active cell A9 ---
$ws.Range("A9:XFD9").Select() | Out-Null
